$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "republishFilter"

$ws.Range("A1").Value = "Editorial Role"
$ws.Range("B1").Value = "PostTitle"
$ws.Range("C1").Value = "Logged_In_Author"
$ws.Range("D1").Value = "AuthorName"
$ws.Range("E1").Value = "Post_type"
$ws.Range("F1").Value = "Buttons"
$ws.Range("G1").Value = "Title_matching "
$ws.Range("H1").Value = "Post_content_match"
$ws.Range("I1").Value = "Title_keyword_match"
$ws.Range("J1").Value = "Post_Keyword_match"
$ws.Range("K1").Value = "Status"

$ws.Columns.Item(1).ColumnWidth = 12.877604166666666
$ws.Columns.Item(2).ColumnWidth = 47.592447916666664
$ws.Columns.Item(3).ColumnWidth = 15.451822916666666
$ws.Columns.Item(4).ColumnWidth = 12.592447916666666
$ws.Columns.Item(5).ColumnWidth = 16.022135416666668
$ws.Columns.Item(6).ColumnWidth = 41.736979166666664
$ws.Columns.Item(7).ColumnWidth = 32.451822916666664
$ws.Columns.Item(8).ColumnWidth = 47.166666666666664
$ws.Columns.Item(9).ColumnWidth = 13.592447916666666
$ws.Columns.Item(10).ColumnWidth = 12.736979166666666

$ws.Range("D9").Select()
